# Generate Report for Handoff
# Renames the source markdown file (new GUID) and the generated xliff
# hand-off artifacts, refreshes the "latest" timestamps, and resets the
# per-language "handback" columns now that a new hand-off round has begun.

$wb = $excel.ActiveWorkbook

$oldGuid = "2c38df77-837e-43d9-b4c1-8276e797efe4"
$newGuid = "c93cea6a-fd96-44c7-9337-f712cf1061ec"
$oldHash = "820350ece49fe2ca1b5fc060f04ed6e44c428dec"
$newHash = "a63a19f9b5b79a25746d0ca06df27f62076054be"

$hyperlinkColor = 15570276  # RGB(100,149,237) "cornflower blue" packed as BGR for OLE_COLOR

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edf86a49926fc6dce97a773cf7af225530c26936/e2e/$oldGuid.md", "", "", "e2e\$newGuid.md")
$wsOverview.Range("B2").Font.Underline = 2
$wsOverview.Range("B2").Font.Color = $hyperlinkColor

$wsOverview.Range("G2").Value = "2016-08-27 04:58:19"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-27 04:58:14"

$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("I2").Font.Underline = -4142
$wsZhCn.Range("I2").Font.ColorIndex = -4105

$wsZhCn.Range("J2").ClearContents()

$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edf86a49926fc6dce97a773cf7af225530c26936/e2e/$oldGuid.md", "", "", "$newGuid.md")
$wsZhCn.Range("A2").Font.Underline = 2
$wsZhCn.Range("A2").Font.Color = $hyperlinkColor

$wsZhCn.Columns.Item(9).ColumnWidth = 17.75
$wsZhCn.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-27 04:58:19"

$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("I2").Font.Underline = -4142
$wsDeDe.Range("I2").Font.ColorIndex = -4105

$wsDeDe.Range("J2").ClearContents()

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edf86a49926fc6dce97a773cf7af225530c26936/e2e/$oldGuid.md", "", "", "$newGuid.md")
$wsDeDe.Range("A2").Font.Underline = 2
$wsDeDe.Range("A2").Font.Color = $hyperlinkColor

$wsDeDe.Columns.Item(9).ColumnWidth = 17.75
$wsDeDe.Columns.Item(10).ColumnWidth = 20.75

Write-Host "Generate Report for Handoff: applied"
